$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item('展览')
$ws.Range('F3').Value = 635
$ws.Range('F4').Value = 2799
$ws.Range('F8').Value = 276
$ws.Range('F9').Value = 6478
$ws.Range('F10').Value = 5
$ws.Range('F13').Value = 5061
$ws.Range('F14').Value = 3
$ws.Range('F16').Value = 2653
$ws.Range('F19').Value = 327
$ws.Range('D21').Value = '澳门路168号 月星国际家居'
$ws.Range('F21').Value = 136
$ws.Range('F22').Value = 1091
$ws.Range('F23').Value = 255
$ws.Range('F25').Value = 1399
$ws.Range('F26').Value = 1066
$ws.Range('F29').Value = 39
$ws.Range('F30').Value = 41
$ws.Range('F31').Value = 112
$ws.Range('F32').Value = 259
$ws.Range('F33').Value = 1526
$ws.Range('F36').Value = 620
$ws.Range('F37').Value = 1085
$ws.Range('F41').Value = 2321
$ws.Range('F42').Value = 2581
$ws.Range('F44').Value = 146

# ---- Sheet: 演出 ----
$ws = $wb.Worksheets.Item('演出')
$ws.Range('F12').Value = 98
$ws.Range('F17').Value = 163
$ws.Range('F19').Value = 2
$ws.Range('F21').Value = 155
$ws.Range('F27').Value = 419

# ---- Sheet: 本地生活 ----
$ws = $wb.Worksheets.Item('本地生活')
$ws.Range('F4').Value = 518
$ws.Range('F8').Value = 1536
$ws.Range('F9').Value = 1823
$ws.Range('F10').Value = 2580
$ws.Range('F11').Value = 885
$ws.Range('F12').Value = 768
$ws.Range('F14').Value = 142

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item('全部类型')
$ws.Range('F3').Value = 518
$ws.Range('F5').Value = 635
$ws.Range('F6').Value = 2799
$ws.Range('F8').Value = 1536
$ws.Range('F9').Value = 276
$ws.Range('F10').Value = 6478
$ws.Range('F11').Value = 885
$ws.Range('F12').Value = 768
$ws.Range('F13').Value = 5061
$ws.Range('F14').Value = 3
$ws.Range('F16').Value = 2653
$ws.Range('F19').Value = 327
$ws.Range('D22').Value = '澳门路168号 月星国际家居'
$ws.Range('F22').Value = 136
$ws.Range('F24').Value = 1091
$ws.Range('F25').Value = 255
$ws.Range('C26').Value = '上海·剑网3×HAPPY ZOO 剑网3十五周年主题咖啡厅'
$ws.Range('D26').Value = '南京东路340号百联zx创趣场四楼05号 HAPPY ZOO'
$ws.Range('E26').Value = '2024.08.16 00:00-10.07 23:59'
$ws.Range('F26').Value = 142
$ws.Range('G26').Value = 10
$ws.Range('H26').Value = 'https://show.bilibili.com/platform/detail.html?id=90305'
$ws.Range('I26').Value = '//i2.hdslb.com/bfs/openplatform/202408/QzCwiYge1722838646403.png'
$ws.Range('B27').Value = '2024-08-17'
$ws.Range('C27').Value = '上海·坏孩纸物语第50届动漫节之豫让篇（免费展）'
$ws.Range('D27').Value = '外马路601号 老码头-3号库'
$ws.Range('E27').Value = '2024.08.17 13:00-08.18 16:00'
$ws.Range('F27').Value = 548
$ws.Range('G27').Value = 36.9
$ws.Range('H27').Value = 'https://show.bilibili.com/platform/detail.html?id=90279'
$ws.Range('I27').Value = '//i0.hdslb.com/bfs/openplatform/202408/oVGVm4Wh1722445973040.png'
$ws.Range('C28').Value = '上海·寻漫岛动漫嘉年华'
$ws.Range('D28').Value = '中山北路3300号4楼L4001号 环球港上海世嘉都市乐园'
$ws.Range('E28').Value = '2024.08.17 10:00-08.18 17:00'
$ws.Range('F28').Value = 1399
$ws.Range('G28').Value = 60
$ws.Range('H28').Value = 'https://show.bilibili.com/platform/detail.html?id=87628'
$ws.Range('I28').Value = '//i1.hdslb.com/bfs/openplatform/202407/ePcJqKzI1721099263380.jpeg'
$ws.Range('C29').Value = '上海·第二届妖妖动漫游戏展'
$ws.Range('D29').Value = '吴中路1588号上海爱琴海购物中心F4 竞梦元宇宙'
$ws.Range('F29').Value = 1066
$ws.Range('G29').Value = 68
$ws.Range('H29').Value = 'https://show.bilibili.com/platform/detail.html?id=90284'
$ws.Range('I29').Value = '//i2.hdslb.com/bfs/openplatform/202408/Q3xelO9p1722578696753.jpeg'
$ws.Range('C30').Value = '上海·第六届燃梦BACG PRO动漫嘉年华-我们在燃梦相遇吧！'
$ws.Range('D30').Value = '盈浦街道淀山浦社区淀山湖大道851号青浦万达茂F3 万达汽车乐园(青浦万达茂店)'
$ws.Range('E30').Value = '2024.08.17 11:00-08.18 16:00'
$ws.Range('F30').Value = 2118
$ws.Range('G30').Value = 65.8
$ws.Range('H30').Value = 'https://show.bilibili.com/platform/detail.html?id=85239'
$ws.Range('I30').Value = '//i1.hdslb.com/bfs/openplatform/202405/mzD4rhY21715109458100.jpeg'
$ws.Range('C31').Value = '上海·魔都特摄同人嘉年华'
$ws.Range('D31').Value = '天山路1111号 现所创邑MIX'
$ws.Range('E31').Value = '2024.08.17 09:30-08.18 17:30'
$ws.Range('F31').Value = 595
$ws.Range('G31').Value = 69
$ws.Range('H31').Value = 'https://show.bilibili.com/platform/detail.html?id=89516'
$ws.Range('I31').Value = '//i0.hdslb.com/bfs/openplatform/202407/0050E5641721292312668.png'
$ws.Range('C32').Value = '上海·（国际）微缩艺术模玩展-GMHS 2024'
$ws.Range('D32').Value = '国展路1099号 上海世博展览馆'
$ws.Range('E32').Value = '2024.08.17 09:00-08.18 17:00'
$ws.Range('F32').Value = 39
$ws.Range('G32').Value = 49
$ws.Range('H32').Value = 'https://show.bilibili.com/platform/detail.html?id=90343'
$ws.Range('I32').Value = '//i2.hdslb.com/bfs/openplatform/202408/L7dY65lZ1722843040285.jpeg'
$ws.Range('B33').Value = '2024-08-18'
$ws.Range('C33').Value = '上海·东方PartyNight'
$ws.Range('D33').Value = '重庆南路308号3楼（近建国中路） 上海MaoLivehouse'
$ws.Range('E33').Value = '2024.08.18 22:00-08.19 01:00'
$ws.Range('F33').Value = 163
$ws.Range('G33').Value = 149
$ws.Range('H33').Value = 'https://show.bilibili.com/platform/detail.html?id=89209'
$ws.Range('I33').Value = '//i1.hdslb.com/bfs/openplatform/202407/4lZtvl551720680564562.jpeg'
$ws.Range('C34').Value = '上海·魅知幻想博览会2024·星辉琉璃 东方同人only'
$ws.Range('D34').Value = '中江路35号 上海跨国采购会展中心'
$ws.Range('E34').Value = '2024.08.18 10:00-08.18 16:00'
$ws.Range('F34').Value = 41
$ws.Range('G34').Value = 88
$ws.Range('H34').Value = 'https://show.bilibili.com/platform/detail.html?id=90366'
$ws.Range('I34').Value = '//i2.hdslb.com/bfs/openplatform/202408/2v7pRUhd1722856549556.jpeg'
$ws.Range('B35').Value = '2024-08-24'
$ws.Range('C35').Value = '上海·BH盛会之星的邀约'
$ws.Range('D35').Value = '鲁班路300号 星光摄影器材城'
$ws.Range('E35').Value = '2024.08.24 10:30-08.25 17:00'
$ws.Range('F35').Value = 259
$ws.Range('G35').Value = 60
$ws.Range('H35').Value = 'https://show.bilibili.com/platform/detail.html?id=88603'
$ws.Range('I35').Value = '//i1.hdslb.com/bfs/openplatform/202407/fUi7Oz2b1719995931315.png'
$ws.Range('C36').Value = '上海·HAG 1st live in Shanghai《不眨眼》2024演唱会'
$ws.Range('D36').Value = '中兴路1683号金融街购物中心三楼L3-27 蜚声LIVE House'
$ws.Range('E36').Value = '2024.08.24 19:30-08.24 21:30'
$ws.Range('F36').Value = 43
$ws.Range('G36').Value = 480
$ws.Range('H36').Value = 'https://show.bilibili.com/platform/detail.html?id=89977'
$ws.Range('I36').Value = '//i1.hdslb.com/bfs/openplatform/202407/iXZNZNM01722243246403.png'
$ws.Range('C37').Value = '上海·coser动漫展'
$ws.Range('D37').Value = '海潮路133号B1 JUMP工坊'
$ws.Range('E37').Value = '2024.08.24 10:00-08.25 17:00'
$ws.Range('F37').Value = 1526
$ws.Range('G37').Value = 60
$ws.Range('H37').Value = 'https://show.bilibili.com/platform/detail.html?id=87347'
$ws.Range('I37').Value = '//i0.hdslb.com/bfs/openplatform/202406/i6vAgX8I1719311206769.jpeg'
$ws.Range('C38').Value = '上海·女仆ONLY同人动漫嘉年华'
$ws.Range('D38').Value = '杨树浦路198号(金茂北外滩)B1层 Terra Park北外滩'
$ws.Range('F38').Value = 7
$ws.Range('G38').Value = 78
$ws.Range('H38').Value = 'https://show.bilibili.com/platform/detail.html?id=90365'
$ws.Range('I38').Value = '//i0.hdslb.com/bfs/openplatform/202408/8NhxowWK1722690320409.jpeg'
$ws.Range('F39').Value = 1085
$ws.Range('F44').Value = 2321
$ws.Range('F45').Value = 2581
$ws.Range('F46').Value = 146

